# Applies the "Add files via upload" edit to libro11.xlsx:
#  - clears the stray "Virtual" value left in H10
#  - appends 10 new patient rows (12-21) with matching text/number formatting
#  - moves the active selection to G25 (as recorded by the author's last save)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108

# --- Row 10: Clinica column no longer populated for this patient ---
$ws.Range("H10").ClearContents() | Out-Null

# --- New patient rows ---
# Columns: A=ID  B=Nombres  C=Apellidos  D=Correo  E=Numero  F=Dias  G=Hora  H=Clinica
$newRows = @(
    @{ Row = 12; A = 1110; B = "Santiago";   C = "Morales Cruz";    D = "santiago.morales99@gmail.com"; E = "7012-1234"; F = "Lunes y miércoles"; G = "9:00 am a 10:00 am";  H = "Virtual" },
    @{ Row = 13; A = 1111; B = "Camila";     C = "Gómez Herrera";   D = "camila.gomez12@yahoo.com";     E = "7721-5678"; F = $null;               G = "1:00 pm a 1:50 pm";  H = "Virtual" },
    @{ Row = 14; A = 1112; B = "Matías";     C = "Castro López";    D = "matias.castro21@hotmail.com";  E = "7534-6789"; F = "Martes";              G = "10:00 am a 12:00 pm"; H = "Virtual" },
    @{ Row = 15; A = 1113; B = "Sofía";      C = "Ramírez Torres";  D = "sofia.ramirez90@gmail.com";    E = "7623-9876"; F = "Jueves";              G = $null;                H = "Virtual" },
    @{ Row = 16; A = 1114; B = "Andrés";     C = "Pérez Martínez";  D = "andres.perez44@gmail.com";     E = "7256-4321"; F = "Viernes";             G = "1:00 pm a 5:00 pm";  H = "Virtual" },
    @{ Row = 17; A = 1115; B = "Isabella";   C = "Hernández Mejía"; D = "isabella.hdz34@gmail.com";     E = "7018-5543"; F = "Lunes o martes";      G = "10:00 am a 11:00 am"; H = "Virtual" },
    @{ Row = 18; A = 1116; B = "Leonardo";   C = "Vargas Díaz";     D = $null;                          E = "7812-2345"; F = "Miércoles";           G = "4:00 pm a 5:00 pm";  H = $null },
    @{ Row = 19; A = 1117; B = "Valentina";  C = "Ortiz Rodríguez"; D = "valentina.ortiz15@gmail.com";  E = "7589-8765"; F = "Viernes";             G = "5:00 pm a 6:00 pm";  H = "Virtual" },
    @{ Row = 20; A = 1118; B = "Martín";     C = "Sánchez Aguilar"; D = "martin.sanchez99@gmail.com";   E = "6842-2233"; F = $null;                 G = "5:00 pm a 6:00 pm";  H = "Virtual" },
    @{ Row = 21; A = 1119; B = "Emilia";     C = "López Morales";   D = "emilia.lopez@hotmail.com";     E = $null;       F = "Lunes";               G = "2:00 pm a 3:00 pm";  H = "Virtual" }
)

# Columns D-13 and D-18 in the source sheet never received the paste-formatting
# pass (they kept the plain column format), so we special-case them below.
$skipFormat = @{ "D13" = $true; "D18" = $true }

foreach ($r in $newRows) {
    $rowNum = $r.Row

    $ws.Rows.Item($rowNum).RowHeight = 14.4

    # Column A: ID number, centered (matches style used in every other row)
    $cellA = $ws.Range("A$rowNum")
    $cellA.Value = $r.A
    $cellA.Style = "Normal"
    $cellA.HorizontalAlignment = $xlCenter

    # Column B: Nombres - vertical-centered + wrap
    $cellB = $ws.Range("B$rowNum")
    if ($null -ne $r.B) { $cellB.Value = $r.B } else { $cellB.ClearContents() | Out-Null }
    $cellB.Style = "Normal"
    $cellB.VerticalAlignment = $xlCenter
    $cellB.WrapText = $true

    # Column C: Apellidos - vertical-centered only
    $cellC = $ws.Range("C$rowNum")
    if ($null -ne $r.C) { $cellC.Value = $r.C } else { $cellC.ClearContents() | Out-Null }
    $cellC.Style = "Normal"
    $cellC.VerticalAlignment = $xlCenter

    # Column D: Correo
    $cellD = $ws.Range("D$rowNum")
    if ($null -ne $r.D) { $cellD.Value = $r.D } else { $cellD.ClearContents() | Out-Null }
    if (-not $skipFormat.ContainsKey("D$rowNum")) {
        $cellD.Style = "Normal"
        $cellD.VerticalAlignment = $xlCenter
        $cellD.WrapText = $true
    }

    # Column E: Numero (kept as text)
    $cellE = $ws.Range("E$rowNum")
    if ($null -ne $r.E) { $cellE.Value = $r.E } else { $cellE.ClearContents() | Out-Null }
    $cellE.Style = "Normal"
    $cellE.VerticalAlignment = $xlCenter
    $cellE.WrapText = $true

    # Column F: Dias
    $cellF = $ws.Range("F$rowNum")
    if ($null -ne $r.F) { $cellF.Value = $r.F } else { $cellF.ClearContents() | Out-Null }
    $cellF.Style = "Normal"
    $cellF.VerticalAlignment = $xlCenter
    $cellF.WrapText = $true

    # Column G: Hora
    $cellG = $ws.Range("G$rowNum")
    if ($null -ne $r.G) { $cellG.Value = $r.G } else { $cellG.ClearContents() | Out-Null }
    $cellG.Style = "Normal"
    $cellG.VerticalAlignment = $xlCenter
    $cellG.WrapText = $true

    # Column H: Clinica
    $cellH = $ws.Range("H$rowNum")
    if ($null -ne $r.H) { $cellH.Value = $r.H } else { $cellH.ClearContents() | Out-Null }
    $cellH.Style = "Normal"
    $cellH.VerticalAlignment = $xlCenter
    $cellH.WrapText = $true
}

# --- Restore the selection recorded in the workbook at last save ---
$ws.Range("G25").Select() | Out-Null
